# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows for "Naranja" (Vega Modelo de Temuco) right
# before the existing row 988, pushing all subsequent rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 988 (old rows 988:1053 shift to 992:1057)
$ws.Rows("988:991").Insert()

# Columns A,B,C,E,F,G,H,I,J are constant for every data row in this sheet.
$constA = 10
$constB = "Vega Modelo de Temuco"
$constC = "La Araucanía"
$constE = 9
$constF = "Fruta"
$constG = 100102
$constH = "Cítricos"
$constI = 100102005
$constJ = "Naranja"

# New row data: Fecha(D), Variedad(K), Calidad(L), Volumen(M), Precio min(N),
# Precio max(O), Precio prom(P), Unidad(Q), Origen(R), Precio/Kg(S), Kg/unidad(T)
$newRows = @(
    @{ Row=988; D=44826; K="Cara cara";  L="Primera"; M=280; N=10000;  O=10000;  P=10000;  Q="`$/bandeja 15 kilos granel"; R="Región de O'Higgins"; S=667; T=15  },
    @{ Row=989; D=44826; K="Cara cara";  L="Tercera"; M=155; N=6000;   O=6000;   P=6000;   Q="`$/bandeja 15 kilos granel"; R="Región de O'Higgins"; S=400; T=15  },
    @{ Row=990; D=44826; K="Navel Late"; L="Especial"; M=250; N=13000; O=13000;  P=13000;  Q="`$/caja 18 kilos granel";    R="Región de O'Higgins"; S=722; T=18  },
    @{ Row=991; D=44826; K="Navel Late"; L="Primera"; M=16;  N=170000; O=180000; P=175000; Q="`$/bins (400 kilos)";        R="Región de O'Higgins"; S=438; T=400 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value2 = $constA
    $ws.Range("B$row").Value2 = $constB
    $ws.Range("C$row").Value2 = $constC
    $ws.Range("D$row").Value2 = $r.D
    $ws.Range("E$row").Value2 = $constE
    $ws.Range("F$row").Value2 = $constF
    $ws.Range("G$row").Value2 = $constG
    $ws.Range("H$row").Value2 = $constH
    $ws.Range("I$row").Value2 = $constI
    $ws.Range("J$row").Value2 = $constJ
    $ws.Range("K$row").Value2 = $r.K
    $ws.Range("L$row").Value2 = $r.L
    $ws.Range("M$row").Value2 = $r.M
    $ws.Range("N$row").Value2 = $r.N
    $ws.Range("O$row").Value2 = $r.O
    $ws.Range("P$row").Value2 = $r.P
    $ws.Range("Q$row").Value2 = $r.Q
    $ws.Range("R$row").Value2 = $r.R
    $ws.Range("S$row").Value2 = $r.S
    $ws.Range("T$row").Value2 = $r.T
}
